$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.231.96"
$ws.Range("E2").Value = "  +0.53%  "

# Row 3
$ws.Range("D3").Value = "2.442.16"
$ws.Range("E3").Value = "  -0.11%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.84"
$ws.Range("E5").Value = "  +0.60%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.60"
$ws.Range("E6").Value = "  +0.22%  "

# Row 7
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
$ws.Range("E8").Value = "  +0.78%  "

# Row 9
$ws.Range("D9").Value = "2.437.56"
$ws.Range("E9").Value = "  -0.45%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.110"
$ws.Range("E10").Value = "  -1.16%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.23"
$ws.Range("E12").Value = "  -1.39%  "

# Row 13
$ws.Range("E13").Value = "  -0.39%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.02"
$ws.Range("E14").Value = "  -0.13%  "

# Row 15
$ws.Range("E15").Value = "  -2.62%  "

# Row 16
$ws.Range("D16").Value = "2.876.27"
$ws.Range("E16").Value = "  +2.75%  "

# Row 17
$ws.Range("D17").Value = "63.019.44"
$ws.Range("E17").Value = "  +0.55%  "

# Row 18
$ws.Range("D18").Value = "2.457.26"
$ws.Range("E18").Value = "  +0.87%  "

# Row 19
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("E20").Value = "  +5.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.96"
$ws.Range("E21").Value = "  +0.85%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.18"
$ws.Range("E22").Value = "  +0.06%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.08"
$ws.Range("E23").Value = "  +11.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.17%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.26"
$ws.Range("E25").Value = "  -3.15%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "613.79"
$ws.Range("E26").Value = "  +3.85%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.84"
$ws.Range("E27").Value = "  +1.39%  "

# Row 28
$ws.Range("E28").Value = "  +0.42%  "

# Row 29
$ws.Range("D29").Value = "2.559.49"
$ws.Range("E29").Value = "  -0.18%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.49"
$ws.Range("E30").Value = "  +2.29%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.26"
$ws.Range("E32").Value = "  -2.53%  "

# Row 33
$ws.Range("E33").Value = "  -4.51%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.89"
$ws.Range("E34").Value = "  +0.61%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.19"
$ws.Range("E35").Value = "  +6.25%  "

# Row 36
$ws.Range("E36").Value = "  -1.59%  "

# Row 37
$ws.Range("E37").Value = "  +0.23%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.379"
$ws.Range("E38").Value = "  -1.02%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.41"
$ws.Range("E39").Value = "  -0.78%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.73"
$ws.Range("E40").Value = "  -0.56%  "

# Row 41
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "146.34"
$ws.Range("E41").Value = "  -1.59%  "

# Row 42
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.69"
$ws.Range("E42").Value = "  +9.18%  "

# Row 43
$ws.Range("E43").Value = "  -2.40%  "

# Row 44
$ws.Range("E44").Value = "  -0.51%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.84"
$ws.Range("E45").Value = "  +0.47%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "148.20"
$ws.Range("E46").Value = "  -0.74%  "

# Row 47
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.75"
$ws.Range("E47").Value = "  +1.27%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.15"
$ws.Range("E48").Value = "  +2.28%  "

# Row 49
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0534"
$ws.Range("E49").Value = "  -0.63%  "

# Row 50
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.599"
$ws.Range("E50").Value = "  -0.66%  "

# Row 51
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0233"
$ws.Range("E51").Value = "  +0.18%  "
